$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new progress note for Suruchi Shrey (row 8) in column B
$ws.Range("B8").Value = "Set up flutter environment and studied basics of dart"

# Reflect the active cell/selection on Sheet1 as edited
$ws.Activate()
$ws.Range("B8").Select()
